$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.709.35"
$ws.Range("E2").Value = "  -5.36%  "
$ws.Range("D3").Value = "3.061.02"
$ws.Range("E3").Value = "  -5.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.053.89"
$ws.Range("E8").Value = "  -5.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("E11").Value = "  -13.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.21%  "
$ws.Range("E14").Value = "  -6.00%  "
$ws.Range("D15").Value = "3.558.55"
$ws.Range("E15").Value = "  -5.82%  "
$ws.Range("D16").Value = "62.675.41"
$ws.Range("E16").Value = "  -5.64%  "
$ws.Range("E17").Value = "  -3.06%  "
$ws.Range("D18").Value = "3.061.51"
$ws.Range("E18").Value = "  -6.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.32%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -14.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.76%  "
$ws.Range("E32").Value = "  -6.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -12.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "490.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -13.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.27%  "
$ws.Range("D38").Value = "3.145.61"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0797"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.72%  "
$ws.Range("E41").Value = "  -10.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -15.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.96%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.65%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("E50").Value = "  -8.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.06%  "
